$wb = $excel.ActiveWorkbook

# --- Sheet "final_fail" (sheet1) updates ---
$ws1 = $wb.Worksheets.Item("final_fail")
$ws1.Cells.Item(2, 1).Value = "Number of clicks"
$ws1.Cells.Item(4, 1).Value = "On/off campus click ratio"
$ws1.Cells.Item(4, 4).Value = $True
$ws1.Cells.Item(4, 10).Value = 8
$ws1.Cells.Item(5, 1).Value = "Days with no interaction (%)"
$ws1.Cells.Item(6, 1).Value = "Largest period of inactivity (h)"
$ws1.Cells.Item(6, 2).Value = $True
$ws1.Cells.Item(6, 3).Value = $True
$ws1.Cells.Item(6, 4).Value = $False
$ws1.Cells.Item(6, 10).Value = 7
$ws1.Cells.Item(7, 1).Value = "Clicks per session"
$ws1.Cells.Item(8, 1).Value = "Clicks per day"
$ws1.Cells.Item(11, 1).Value = "Resources viewed"
$ws1.Cells.Item(12, 1).Value = "Average session duration (min)"
$ws1.Cells.Item(13, 1).Value = "Start of Session 2 (%)"
$ws1.Cells.Item(14, 1).Value = "Start of Session 1 (%)"
$ws1.Cells.Item(15, 1).Value = "Days with no interaction"
$ws1.Cells.Item(16, 1).Value = "Clicks on course"
$ws1.Cells.Item(17, 1).Value = "Start of Session 4 (%)"
$ws1.Cells.Item(17, 3).Value = $True
$ws1.Cells.Item(17, 5).Value = $False
$ws1.Cells.Item(18, 1).Value = "Clicks on campus"
$ws1.Cells.Item(18, 2).Value = $False
$ws1.Cells.Item(18, 5).Value = $True
$ws1.Cells.Item(19, 2).Value = $False
$ws1.Cells.Item(19, 10).Value = 2
$ws1.Cells.Item(20, 1).Value = "Files downloaded"
$ws1.Cells.Item(21, 1).Value = "Quizzes started"
$ws1.Cells.Item(22, 1).Value = "Assignments viewed"
$ws1.Cells.Item(23, 1).Value = "Clicks on forum"
$ws1.Cells.Item(23, 2).Value = $True
$ws1.Cells.Item(23, 4).Value = $False
$ws1.Cells.Item(24, 1).Value = "Clicks on folder"
$ws1.Cells.Item(24, 4).Value = $True
$ws1.Cells.Item(24, 10).Value = 2
$ws1.Cells.Item(25, 1).Value = "Assignments submitted"
$ws1.Cells.Item(26, 1).Value = "Discussions viewed"
$ws1.Cells.Item(27, 1).Value = "Forum posts"
$ws1.Cells.Item(28, 1).Value = "Number of sessions"
$ws1.Cells.Item(29, 1).Value = "Start of Session 10 (%)"
$ws1.Cells.Item(30, 1).Value = "Number of days"

# --- Sheet "final_gifted" (sheet2) updates ---
$ws2 = $wb.Worksheets.Item("final_gifted")
$ws2.Cells.Item(6, 1).Value = "Average session duration (min)"
$ws2.Cells.Item(7, 1).Value = "Clicks per day"
$ws2.Cells.Item(8, 1).Value = "Clicks per session"
$ws2.Cells.Item(9, 1).Value = "Days with no interaction"
$ws2.Cells.Item(9, 4).Value = $False
$ws2.Cells.Item(9, 10).Value = 5
$ws2.Cells.Item(10, 1).Value = "On/off campus click ratio"
$ws2.Cells.Item(16, 1).Value = "Clicks on campus"
$ws2.Cells.Item(16, 2).Value = $False
$ws2.Cells.Item(16, 4).Value = $True
$ws2.Cells.Item(17, 1).Value = "Resources viewed"
$ws2.Cells.Item(17, 2).Value = $True
$ws2.Cells.Item(17, 10).Value = 5
$ws2.Cells.Item(26, 1).Value = "Start of Session 10 (%)"
$ws2.Cells.Item(27, 1).Value = "Forum posts"
$ws2.Cells.Item(27, 3).Value = $False
$ws2.Cells.Item(27, 4).Value = $True
$ws2.Cells.Item(28, 1).Value = "Start of Session 9 (%)"
$ws2.Cells.Item(28, 3).Value = $True
$ws2.Cells.Item(28, 4).Value = $False
$ws2.Cells.Item(36, 1).Value = "Discussions viewed"
